$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of EUR->ARS rate history data (row 45)
# Force column A to be stored as text so the date-like string
# "2025-09-28" is not auto-converted into a date serial number,
# then clear the temporary formatting so no extra style is left
# on the cell.
$ws.Range("A45").NumberFormat = "@"
$ws.Range("A45").Value = "2025-09-28"
$ws.Range("A45").ClearFormats()

$ws.Range("B45").Value = "15:16:37"
$ws.Range("C45").Value = "1.00 EUR = 1,623.5666"
